$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = "26.110.59"
$ws.Cells.Item(2,5).Value = "  -0.86%  "
$ws.Cells.Item(3,4).Value = "1.652.48"
$ws.Cells.Item(3,5).Value = "  -1.00%  "
$ws.Cells.Item(4,5).Value = "  -0.31%  "
$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value = "218.83"
$ws.Cells.Item(5,5).Value = "  -0.81%  "
$ws.Cells.Item(6,4).NumberFormat = "@"
$ws.Cells.Item(6,4).Value = "0.5240"
$ws.Cells.Item(6,5).Value = "  -1.42%  "
$ws.Cells.Item(7,5).Value = "  -0.31%  "
$ws.Cells.Item(8,4).NumberFormat = "@"
$ws.Cells.Item(8,4).Value = "0.2659"
$ws.Cells.Item(8,5).Value = "  +0.28%  "
$ws.Cells.Item(9,4).NumberFormat = "@"
$ws.Cells.Item(9,4).Value = "0.06341"
$ws.Cells.Item(9,5).Value = "  -0.44%  "
$ws.Cells.Item(10,4).NumberFormat = "@"
$ws.Cells.Item(10,4).Value = "20.63"
$ws.Cells.Item(10,5).Value = "  -1.57%  "
$ws.Cells.Item(11,4).NumberFormat = "@"
$ws.Cells.Item(11,4).Value = "0.07699"
$ws.Cells.Item(11,5).Value = "  -1.80%  "
$ws.Cells.Item(12,4).NumberFormat = "@"
$ws.Cells.Item(12,4).Value = "4.612"
$ws.Cells.Item(12,5).Value = "  +1.74%  "
$ws.Cells.Item(13,4).Value = "1.576.73"
$ws.Cells.Item(13,5).Value = "  -5.78%  "
$ws.Cells.Item(14,4).Value = "1.879.59"
$ws.Cells.Item(14,5).Value = "  -0.96%  "
$ws.Cells.Item(15,4).NumberFormat = "@"
$ws.Cells.Item(15,4).Value = "0.5609"
$ws.Cells.Item(15,5).Value = "  -0.21%  "
$ws.Cells.Item(16,4).Value = "0.0₅8201"
$ws.Cells.Item(16,5).Value = "  +0.54%  "
$ws.Cells.Item(17,4).NumberFormat = "@"
$ws.Cells.Item(17,4).Value = "65.26"
$ws.Cells.Item(17,5).Value = "  -0.97%  "
$ws.Cells.Item(18,4).Value = "26.121.08"
$ws.Cells.Item(18,5).Value = "  -0.79%  "
$ws.Cells.Item(19,4).NumberFormat = "@"
$ws.Cells.Item(19,4).Value = "1.004"
$ws.Cells.Item(19,5).Value = "  -0.32%  "
$ws.Cells.Item(20,4).NumberFormat = "@"
$ws.Cells.Item(20,4).Value = "4.700"
$ws.Cells.Item(20,5).Value = "  -0.50%  "
$ws.Cells.Item(21,4).NumberFormat = "@"
$ws.Cells.Item(21,4).Value = "10.37"
$ws.Cells.Item(21,5).Value = "  +0.86%  "
$ws.Cells.Item(22,4).NumberFormat = "@"
$ws.Cells.Item(22,4).Value = "191.04"
$ws.Cells.Item(22,5).Value = "  -3.78%  "
$ws.Cells.Item(23,4).NumberFormat = "@"
$ws.Cells.Item(23,4).Value = "5.984"
$ws.Cells.Item(23,5).Value = "  -1.26%  "
$ws.Cells.Item(24,4).NumberFormat = "@"
$ws.Cells.Item(24,4).Value = "1.005"
$ws.Cells.Item(24,5).Value = "  -0.33%  "
$ws.Cells.Item(25,4).NumberFormat = "@"
$ws.Cells.Item(25,4).Value = "145.32"
$ws.Cells.Item(25,5).Value = "  -0.85%  "
$ws.Cells.Item(26,5).Value = "  -1.30%  "
$ws.Cells.Item(27,4).NumberFormat = "@"
$ws.Cells.Item(27,4).Value = "7.257"
$ws.Cells.Item(27,5).Value = "  +0.03%  "
$ws.Cells.Item(28,5).Value = "  -1.52%  "
$ws.Cells.Item(29,4).NumberFormat = "@"
$ws.Cells.Item(29,4).Value = "1.505"
$ws.Cells.Item(29,5).Value = "  -0.17%  "
$ws.Cells.Item(30,4).NumberFormat = "@"
$ws.Cells.Item(30,4).Value = "0.05626"
$ws.Cells.Item(30,5).Value = "  -4.80%  "
$ws.Cells.Item(31,4).NumberFormat = "@"
$ws.Cells.Item(31,4).Value = "1.274"
$ws.Cells.Item(31,5).Value = "  -0.90%  "
$ws.Cells.Item(32,4).NumberFormat = "@"
$ws.Cells.Item(32,4).Value = "3.493"
$ws.Cells.Item(32,5).Value = "  -1.72%  "
$ws.Cells.Item(33,4).NumberFormat = "@"
$ws.Cells.Item(33,4).Value = "3.373"
$ws.Cells.Item(33,5).Value = "  +1.36%  "
$ws.Cells.Item(34,4).NumberFormat = "@"
$ws.Cells.Item(34,4).Value = "1.580"
$ws.Cells.Item(34,5).Value = "  -1.63%  "
$ws.Cells.Item(35,5).Value = "  -1.26%  "
$ws.Cells.Item(36,4).NumberFormat = "@"
$ws.Cells.Item(36,4).Value = "0.9504"
$ws.Cells.Item(36,5).Value = "  -1.32%  "
$ws.Cells.Item(37,5).Value = "  -0.91%  "
$ws.Cells.Item(38,4).NumberFormat = "@"
$ws.Cells.Item(38,4).Value = "0.5754"
$ws.Cells.Item(38,5).Value = "  -1.17%  "
$ws.Cells.Item(39,4).NumberFormat = "@"
$ws.Cells.Item(39,4).Value = "0.01594"
$ws.Cells.Item(39,5).Value = "  -1.48%  "
$ws.Cells.Item(40,4).NumberFormat = "@"
$ws.Cells.Item(40,4).Value = "5.991"
$ws.Cells.Item(40,5).Value = "  +0.56%  "
$ws.Cells.Item(41,5).Value = "  -0.36%  "
$ws.Cells.Item(42,4).NumberFormat = "@"
$ws.Cells.Item(42,4).Value = "0.8387"
$ws.Cells.Item(42,5).Value = "  -2.29%  "
$ws.Cells.Item(43,4).NumberFormat = "@"
$ws.Cells.Item(43,4).Value = "101.74"
$ws.Cells.Item(43,5).Value = "  -1.14%  "
$ws.Cells.Item(44,4).Value = "1.014.05"
$ws.Cells.Item(44,5).Value = "  -5.82%  "
$ws.Cells.Item(45,4).Value = "1.790.45"
$ws.Cells.Item(45,5).Value = "  -0.97%  "
$ws.Cells.Item(46,4).NumberFormat = "@"
$ws.Cells.Item(46,4).Value = "58.28"
$ws.Cells.Item(46,5).Value = "  -0.63%  "
$ws.Cells.Item(47,2).Value = "Frax"
$ws.Cells.Item(47,3).Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Cells.Item(47,4).NumberFormat = "@"
$ws.Cells.Item(47,4).Value = "1.006"
$ws.Cells.Item(47,5).Value = "  -0.80%  "
$ws.Cells.Item(48,2).Value = "Cronos"
$ws.Cells.Item(48,3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(48,4).NumberFormat = "@"
$ws.Cells.Item(48,4).Value = "0.05329"
$ws.Cells.Item(48,5).Value = "  +3.45%  "
$ws.Cells.Item(49,2).Value = "Mantle"
$ws.Cells.Item(49,3).Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Cells.Item(49,4).NumberFormat = "@"
$ws.Cells.Item(49,4).Value = "0.4345"
$ws.Cells.Item(49,5).Value = "  -1.31%  "
$ws.Cells.Item(50,2).Value = "EnergySwap"
$ws.Cells.Item(50,3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(50,4).NumberFormat = "@"
$ws.Cells.Item(50,4).Value = "8.005"
$ws.Cells.Item(50,5).Value = "  -0.58%  "
$ws.Cells.Item(51,2).Value = "Algorand"
$ws.Cells.Item(51,3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(51,4).NumberFormat = "@"
$ws.Cells.Item(51,4).Value = "0.09760"
$ws.Cells.Item(51,5).Value = "  +1.40%  "
